# Commit: "simplifying units with pint and UnitAttribute"
#
# The shared unit-abbreviation strings used throughout the workbook (in the
# "Units" columns of many sheets) are being replaced by their
# pint-compatible long-form spellings. Since these abbreviations are reused
# by many cells across many sheets, we do a whole-cell-value Find & Replace
# on every worksheet for each abbreviation -> long form pair. This updates
# every cell that held the old value, without touching cells that merely
# contain the abbreviation as a substring (e.g. "s" must not touch "s^-1").

$wb = $excel.ActiveWorkbook

$replacements = @(
    @("s", "second"),
    @("C", "degC"),
    @("g", "gram"),
    @("l", "liter"),
    @("M", "molar"),
    @("s^-1", "1 / second"),
    @("M s^-1", "molar / second"),
    @("molecule mol^-1", "molecule / mole"),
    @("mmol/gDCW/h", "millimole / gDCW / hour"),
    @("g l^-1", "gram / liter"),
    @("molecule^-1 s^-1", "1 / molecule / second"),
    @("molecule^-2 s^-1", "1 / molecule ** 2 / second")
)

foreach ($ws in $wb.Worksheets) {
    foreach ($pair in $replacements) {
        # LookAt:=xlWhole (1), MatchCase:=True -- only swap cells whose
        # entire value is the exact (case-sensitive) abbreviation.
        $ws.Cells.Replace($pair[0], $pair[1], 1, 1, $true)
    }
}

# Window/view tweaks recorded in the commit.
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).Height = 13965

$modelSheet = $wb.Worksheets.Item("Model")
$modelSheet.Activate()
$modelSheet.Range("A1").Select()
